$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 2179422
$ws.Range("I33").Value = 2469952.2
$ws.Range("J33").Value = 445
$ws.Range("K33").Value = 2469952.2
$ws.Range("L33").Value = 445
$ws.Range("M33").Value = -2469723.2
$ws.Range("N33").Value = -903
$ws.Range("H88").Value = 2928.138
$ws.Range("I88").Value = 1643.5
$ws.Range("J88").Value = 3023.2964
$ws.Range("K88").Value = 1643.5
$ws.Range("L88").Value = 3023.2964
$ws.Range("M88").Value = -1237.5
$ws.Range("N88").Value = -3835.2964
$ws.Range("H91").Value = 2928.138
$ws.Range("I91").Value = 1643.5
$ws.Range("J91").Value = 3023.2964
$ws.Range("K91").Value = 1643.5
$ws.Range("L91").Value = 3023.2964
$ws.Range("M91").Value = -239.5
$ws.Range("N91").Value = -5831.2964
$ws.Range("H100").Value = 2319.3635
$ws.Range("I100").Value = 2657.111
$ws.Range("J100").Value = 799.5
$ws.Range("K100").Value = 2657.111
$ws.Range("L100").Value = 799.5
$ws.Range("M100").Value = -2116.111
$ws.Range("N100").Value = -1881.5
$ws.Range("H116").Value = 5302.8184
$ws.Range("I116").Value = 4463.3335
$ws.Range("J116").Value = 5617.625
$ws.Range("K116").Value = 4463.3335
$ws.Range("L116").Value = 5617.625
$ws.Range("M116").Value = -1021.3335
$ws.Range("N116").Value = -12501.625
$ws.Range("H132").Value = 22730338
$ws.Range("I132").Value = 23812640
$ws.Range("K132").Value = 71437920
$ws.Range("M132").Value = -71435390
$ws.Range("H135").Value = 1666.3928
$ws.Range("I135").Value = 719.7368
$ws.Range("K135").Value = 6477.6312
$ws.Range("M135").Value = -3942.6312
$ws.Range("H138").Value = 2015.34
$ws.Range("I138").Value = 1085.8889
$ws.Range("K138").Value = 3257.6667
$ws.Range("M138").Value = 1882.3333
$ws.Range("H141").Value = 1690.4706
$ws.Range("I141").Value = 1168.909
$ws.Range("K141").Value = 3506.727
$ws.Range("M141").Value = 1673.273

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4455.5864
$ws.Range("I32").Value = 2615.7637
$ws.Range("K32").Value = 2615.7637
$ws.Range("M32").Value = -2328.7637
$ws.Range("H45").Value = 11069764
$ws.Range("I45").Value = 20552248
$ws.Range("J45").Value = 6865.5
$ws.Range("K45").Value = 20552248
$ws.Range("L45").Value = 6865.5
$ws.Range("M45").Value = -20551871
$ws.Range("N45").Value = -7619.5
$ws.Range("H61").Value = 2849.5186
$ws.Range("I61").Value = 2667.2727
$ws.Range("J61").Value = 3651.4
$ws.Range("K61").Value = 2667.2727
$ws.Range("L61").Value = 3651.4
$ws.Range("M61").Value = -2455.2727
$ws.Range("N61").Value = -4075.4
$ws.Range("H74").Value = 54423.945
$ws.Range("I74").Value = 6947.3413
$ws.Range("J74").Value = 204157.84
$ws.Range("K74").Value = 6947.3413
$ws.Range("L74").Value = 204157.84
$ws.Range("M74").Value = -6073.3413
$ws.Range("N74").Value = -205905.84
$ws.Range("H77").Value = 54423.945
$ws.Range("I77").Value = 6947.3413
$ws.Range("J77").Value = 204157.84
$ws.Range("K77").Value = 34736.7065
$ws.Range("L77").Value = 1020789.2
$ws.Range("M77").Value = -30368.7065
$ws.Range("N77").Value = -1029525.2
$ws.Range("H97").Value = 1085052.2
$ws.Range("I97").Value = 1474760.9
$ws.Range("K97").Value = 1474760.9
$ws.Range("M97").Value = -1474264.9
$ws.Range("H102").Value = 1853993.1
$ws.Range("I102").Value = 2194838
$ws.Range("K102").Value = 2194838
$ws.Range("M102").Value = -2193216
$ws.Range("H124").Value = 12187.167
$ws.Range("J124").Value = 12187.167
$ws.Range("L124").Value = 12187.167
$ws.Range("N124").Value = -22007.167
$ws.Range("H132").Value = 1571.4865
$ws.Range("I132").Value = 1208.2727
$ws.Range("K132").Value = 3624.8181
$ws.Range("M132").Value = -1094.8181
$ws.Range("H136").Value = 2849.5186
$ws.Range("I136").Value = 2667.2727
$ws.Range("J136").Value = 3651.4
$ws.Range("K136").Value = 8001.8181
$ws.Range("L136").Value = 10954.2
$ws.Range("M136").Value = -5451.8181
$ws.Range("N136").Value = -16054.2
$ws.Range("H139").Value = 91660
$ws.Range("J139").Value = 91660
$ws.Range("L139").Value = 91660
$ws.Range("N139").Value = -101940

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 15208.556
$ws.Range("J81").Value = 15208.556
$ws.Range("L81").Value = 15208.556
$ws.Range("N81").Value = -17330.556
$ws.Range("H84").Value = 15208.556
$ws.Range("J84").Value = 15208.556
$ws.Range("L84").Value = 45625.66800000001
$ws.Range("N84").Value = -56233.66800000001
$ws.Range("H105").Value = 2502189
$ws.Range("I105").Value = 3474499.5
$ws.Range("K105").Value = 3474499.5
$ws.Range("M105").Value = -3472752.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 199.88889
$ws.Range("I7").Value = 128.14285
$ws.Range("K7").Value = 128.14285
$ws.Range("M7").Value = -15.14285000000001
$ws.Range("H31").Value = 48528.6
$ws.Range("J31").Value = 68477.07000000001
$ws.Range("L31").Value = 68477.07000000001
$ws.Range("N31").Value = -69067.07000000001
$ws.Range("H34").Value = 48528.6
$ws.Range("J34").Value = 68477.07000000001
$ws.Range("L34").Value = 68477.07000000001
$ws.Range("N34").Value = -68881.07000000001
$ws.Range("H105").Value = 1533.6
$ws.Range("I105").Value = 1643.5834
$ws.Range("J105").Value = 1093.6666
$ws.Range("K105").Value = 1643.5834
$ws.Range("L105").Value = 1093.6666
$ws.Range("M105").Value = 103.4166
$ws.Range("N105").Value = -4587.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 601.75
$ws.Range("I3").Value = 655
$ws.Range("J3").Value = 229
$ws.Range("K3").Value = 1965
$ws.Range("L3").Value = 687
$ws.Range("M3").Value = -1853
$ws.Range("N3").Value = -911
$ws.Range("H5").Value = 1297.7142
$ws.Range("I5").Value = 934.125
$ws.Range("J5").Value = 1782.5
$ws.Range("K5").Value = 2802.375
$ws.Range("L5").Value = 5347.5
$ws.Range("M5").Value = -2690.375
$ws.Range("N5").Value = -5571.5
$ws.Range("H120").Value = 14580.286
$ws.Range("I120").Value = 6343
$ws.Range("K120").Value = 19029
$ws.Range("M120").Value = -14191
$ws.Range("H132").Value = 2457
$ws.Range("J132").Value = 2606.2222
$ws.Range("L132").Value = 23455.9998
$ws.Range("N132").Value = -28515.9998
$ws.Range("H133").Value = 1900.5
$ws.Range("I133").Value = 1900.5
$ws.Range("K133").Value = 5701.5
$ws.Range("M133").Value = -641.5
$ws.Range("H134").Value = 3703
$ws.Range("I134").Value = 3703
$ws.Range("K134").Value = 11109
$ws.Range("M134").Value = -6039
$ws.Range("H135").Value = 1297.7142
$ws.Range("I135").Value = 934.125
$ws.Range("J135").Value = 1782.5
$ws.Range("K135").Value = 8407.125
$ws.Range("L135").Value = 16042.5
$ws.Range("M135").Value = -5872.125
$ws.Range("N135").Value = -21112.5
$ws.Range("H139").Value = 2659.5881
$ws.Range("I139").Value = 2485.5
$ws.Range("K139").Value = 7456.5
$ws.Range("M139").Value = -2316.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5888008.5
$ws.Range("I70").Value = 8005495.5
$ws.Range("K70").Value = 8005495.5
$ws.Range("M70").Value = -8005225.5
$ws.Range("H73").Value = 5888008.5
$ws.Range("I73").Value = 8005495.5
$ws.Range("K73").Value = 8005495.5
$ws.Range("M73").Value = -8004559.5
$ws.Range("H80").Value = 2713990.8
$ws.Range("I80").Value = 6098978.5
$ws.Range("J80").Value = 6000.4
$ws.Range("K80").Value = 6098978.5
$ws.Range("L80").Value = 6000.4
$ws.Range("M80").Value = -6097980.5
$ws.Range("N80").Value = -7996.4
$ws.Range("H83").Value = 2713990.8
$ws.Range("I83").Value = 6098978.5
$ws.Range("J83").Value = 6000.4
$ws.Range("K83").Value = 30494892.5
$ws.Range("L83").Value = 30002
$ws.Range("M83").Value = -30489900.5
$ws.Range("N83").Value = -39986

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H112").Value = 10387
$ws.Range("J112").Value = 10387
$ws.Range("L112").Value = 10387
$ws.Range("N112").Value = -13341
$ws.Range("H132").Value = 6420.8374
$ws.Range("I132").Value = 7029.207
$ws.Range("J132").Value = 5160.643
$ws.Range("K132").Value = 21087.621
$ws.Range("L132").Value = 15481.929
$ws.Range("M132").Value = -18557.621
$ws.Range("N132").Value = -20541.929
$ws.Range("H136").Value = 26013.932
$ws.Range("I136").Value = 34267.938
$ws.Range("K136").Value = 102803.814
$ws.Range("M136").Value = -100253.814

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 18538678
$ws.Range("I132").Value = 23813498
$ws.Range("K132").Value = 71440494
$ws.Range("M132").Value = -71437964
$ws.Range("H136").Value = 3149.3948
$ws.Range("I136").Value = 2660.2
$ws.Range("K136").Value = 7980.599999999999
$ws.Range("M136").Value = -5430.599999999999
